$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 160.9723563333334
$ws.Range("H2").Value = 482.917069
$ws.Range("I2").Value = 0.3931645655589854
$ws.Range("J2").Value = 0.3931645655589854
$ws.Range("M2").Value = 3.759736666666667
$ws.Range("N2").Value = 11.27921
$ws.Range("O2").Value = 0.0683751702595819
$ws.Range("P2").Value = 0.06837517025958188
$ws.Range("Q2").Value = 605.2136704261655
$ws.Range("R2").Value = 5446.92303383549
$ws.Range("S2").Value = 0.02688269411013017
$ws.Range("T2").Value = 0.02688269411013017
$ws.Range("G3").Value = 160.9723563333334
$ws.Range("H3").Value = 482.917069
$ws.Range("I3").Value = 0.3931645655589854
$ws.Range("J3").Value = 0.3931645655589854
$ws.Range("O3").Value = 0.6514180024294648
$ws.Range("P3").Value = 0.6514180024294647
$ws.Range("Q3").Value = 5765.939283738289
$ws.Range("R3").Value = 51893.4535536446
$ws.Range("S3").Value = 0.2561144759224826
$ws.Range("T3").Value = 0.2561144759224825
$ws.Range("G4").Value = 160.9723563333334
$ws.Range("H4").Value = 482.917069
$ws.Range("I4").Value = 0.3931645655589854
$ws.Range("J4").Value = 0.3931645655589854
$ws.Range("O4").Value = 0.2802068273109533
$ws.Range("P4").Value = 0.2802068273109533
$ws.Range("Q4").Value = 2480.213238102579
$ws.Range("R4").Value = 22321.91914292321
$ws.Range("S4").Value = 0.1101673955263726
$ws.Range("T4").Value = 0.1101673955263726
$ws.Range("I5").Value = 0.2197635343237224
$ws.Range("J5").Value = 0.2197635343237224
$ws.Range("M5").Value = 3.759736666666667
$ws.Range("N5").Value = 11.27921
$ws.Range("O5").Value = 0.0683751702595819
$ws.Range("P5").Value = 0.06837517025958188
$ws.Range("Q5").Value = 338.2906469325055
$ws.Range("R5").Value = 3044.615822392549
$ws.Range("S5").Value = 0.01502636907623199
$ws.Range("T5").Value = 0.01502636907623199
$ws.Range("I6").Value = 0.2197635343237224
$ws.Range("J6").Value = 0.2197635343237224
$ws.Range("O6").Value = 0.6514180024294648
$ws.Range("P6").Value = 0.6514180024294647
$ws.Range("S6").Value = 0.1431579225359983
$ws.Range("T6").Value = 0.1431579225359983
$ws.Range("I7").Value = 0.2197635343237224
$ws.Range("J7").Value = 0.2197635343237224
$ws.Range("O7").Value = 0.2802068273109533
$ws.Range("P7").Value = 0.2802068273109533
$ws.Range("S7").Value = 0.06157924271149203
$ws.Range("T7").Value = 0.06157924271149203
$ws.Range("I8").Value = 0.3870719001172923
$ws.Range("J8").Value = 0.3870719001172923
$ws.Range("M8").Value = 3.759736666666667
$ws.Range("N8").Value = 11.27921
$ws.Range("O8").Value = 0.0683751702595819
$ws.Range("P8").Value = 0.06837517025958188
$ws.Range("Q8").Value = 595.8349910189734
$ws.Range("R8").Value = 5362.51491917076
$ws.Range("S8").Value = 0.02646610707321973
$ws.Range("T8").Value = 0.02646610707321973
$ws.Range("I9").Value = 0.3870719001172923
$ws.Range("J9").Value = 0.3870719001172923
$ws.Range("O9").Value = 0.6514180024294648
$ws.Range("P9").Value = 0.6514180024294647
$ws.Range("Q9").Value = 5676.587541261226
$ws.Range("R9").Value = 51089.28787135103
$ws.Range("S9").Value = 0.2521456039709838
$ws.Range("T9").Value = 0.2521456039709838
$ws.Range("I10").Value = 0.3870719001172923
$ws.Range("J10").Value = 0.3870719001172923
$ws.Range("O10").Value = 0.2802068273109533
$ws.Range("P10").Value = 0.2802068273109533
$ws.Range("S10").Value = 0.1084601890730887
$ws.Range("T10").Value = 0.1084601890730887
